$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "comprobado" column (H) -----------------------------------------
# Header in H1
$ws.Range("H1").Value = "comprobado"

# Make sure every data row (2-101) gets the same visual style as the
# handful of H cells that already existed in the sheet (style carries a red
# font, borders, centered alignment -- same "s=5" cellXf used elsewhere in
# the sheet). Cell H13 already uses that exact style, so copy its format
# across the whole H column first.
$ws.Range("H13").Copy() | Out-Null
$ws.Range("H1:H101").PasteSpecial(-4122) | Out-Null

# Now stamp the literal text "true" into H2:H101. A bare Value/Formula
# assignment of the word "true" is auto-coerced to a Boolean by Excel, so
# instead write a formula that evaluates to the text string "true" and
# then convert the whole range to static values (xlPasteValues) -- this
# keeps the cell's stored type as shared-string text instead of boolean.
$data = $ws.Range("H2:H101")
$data.Formula = '="true"'
$data.Copy() | Out-Null
$data.PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# --- Selection --------------------------------------------------------
$ws.Range("I11").Select() | Out-Null
